$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1 - copy the style of the existing header cell H1
# so they match the other header cells (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-21: column I (I0) and column J (IF) values
$data = @(
    @(1, 2),
    @(1, 4),
    @(4, 6),
    @(6, 8),
    @(3, 5),
    @(7, 8),
    @(9, 9),
    @(7, 9),
    @(3, 4),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(5, 8),
    @(1, 4),
    @(4, 6),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
